$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range keeps its original text formatting so that purely
# numeric-looking values (e.g. "10.17") are not auto-converted into numbers.
$ws.Range('B2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '57.477.97'
$ws.Range('E2').Value = '  -2.30%  '
$ws.Range('D3').Value = '2.566.10'
$ws.Range('E3').Value = '  -3.42%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '520.54'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = '143.66'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  -1.28%  '
$ws.Range('D9').Value = '2.583.43'
$ws.Range('E10').Value = '  -5.37%  '
$ws.Range('E11').Value = '  -2.41%  '
$ws.Range('D12').Value = '0.325'
$ws.Range('E12').Value = '  -2.90%  '
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = '3.022.45'
$ws.Range('E14').Value = '  -3.26%  '
$ws.Range('D15').Value = '57.484.67'
$ws.Range('E15').Value = '  -2.25%  '
$ws.Range('D16').Value = '20.14'
$ws.Range('E16').Value = '  -4.04%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.589.92'
$ws.Range('E17').Value = '  -2.79%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.0000133'
$ws.Range('E18').Value = '  -2.52%  '
$ws.Range('D19').Value = '334.97'
$ws.Range('E19').Value = '  -1.21%  '
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('D21').Value = '10.17'
$ws.Range('E21').Value = '  -2.03%  '
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').Value = '64.58'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = '0.166'
$ws.Range('E25').Value = '  -0.19%  '
$ws.Range('D26').Value = '0.401'
$ws.Range('E26').Value = '  -4.55%  '
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').Value = '2.676.46'
$ws.Range('E28').Value = '  -3.69%  '
$ws.Range('D29').Value = '6.95'
$ws.Range('E29').Value = '  -2.80%  '
$ws.Range('D30').Value = '0.0₃0749'
$ws.Range('E30').Value = '  -6.55%  '
$ws.Range('E31').Value = '  -0.10%  '
$ws.Range('D32').Value = '6.23'
$ws.Range('E32').Value = '  -6.58%  '
$ws.Range('E33').Value = '  -0.60%  '
$ws.Range('D34').Value = '18.58'
$ws.Range('E34').Value = '  -1.61%  '
$ws.Range('D35').Value = '148.43'
$ws.Range('E35').Value = '  -1.48%  '
$ws.Range('D36').Value = '4.02'
$ws.Range('E36').Value = '  -3.05%  '
$ws.Range('E37').Value = '  -4.28%  '
$ws.Range('D38').Value = '0.837'
$ws.Range('E38').Value = '  -9.49%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').Value = '35.84'
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '0.833'
$ws.Range('E40').Value = '  -4.08%  '
$ws.Range('E41').Value = '  -1.54%  '
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '267.58'
$ws.Range('E44').Value = '  -2.77%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = '10.65'
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.0953'
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.588'
$ws.Range('E47').Value = '  -3.97%  '
$ws.Range('D48').Value = '18.88'
$ws.Range('E48').Value = '  -4.28%  '
$ws.Range('D49').Value = '0.0520'
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '1.970.61'
$ws.Range('E50').Value = '  -4.28%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '4.56'
$ws.Range('E51').Value = '  -3.55%  '
